$d = $word.ActiveDocument

# --------------------------------------------------------------------
# 0) The document starts out with a "_GoBack" bookmark collapsed right
#    after "...estadisticas personaje." (the very last paragraph). The
#    edit relocates that bookmark so it wraps the "PJ_active" item
#    instead, so remove the old one first.
# --------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete() | Out-Null
}

# --------------------------------------------------------------------
# 1) "Revisar que funcione correctamente la variable PJ_active."
#    -> strike through the whole paragraph (runs + paragraph mark) and
#       wrap the visible text with the relocated "_GoBack" bookmark.
# --------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute(
    "Revisar que funcione correctamente la variable PJ_active.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $para1 = $rng1.Paragraphs(1)

    $full1 = $para1.Range
    $full1.Font.StrikeThrough = 1

    $textOnly1 = $para1.Range
    $textOnly1.MoveEnd(1, -1) | Out-Null
    $d.Bookmarks.Add("_GoBack", $textOnly1) | Out-Null
}

# --------------------------------------------------------------------
# 2) "Cuando se carga city.html ... modifique los datos?)"
#    -> strike through the whole paragraph (runs + paragraph mark).
# --------------------------------------------------------------------
$rng2 = $d.Content
$needle2 = "Cuando se carga city.html hay momentos donde toma los datos de " + `
    "PJ_active y otros donde dice que es null, revisar. (ver si se puede " + `
    "colocar que cargue todos los datos antes o que cargue la pagina y " + `
    "despu" + [char]0xE9 + "s modifique los datos?)"
$found2 = $rng2.Find.Execute(
    $needle2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $para2 = $rng2.Paragraphs(1)
    $full2 = $para2.Range
    $full2.Font.StrikeThrough = 1
}

# --------------------------------------------------------------------
# 3) "Revisar calculos de estadisticas personaje."
#    -> strike through everything except the trailing period, which
#       becomes its own, unstruck run. (The old "_GoBack" bookmark
#       that used to sit here was already removed in step 0.)
# --------------------------------------------------------------------
$rng3 = $d.Content
$needle3 = "Revisar c" + [char]0xE1 + "lculos de estad" + [char]0xED + "sticas personaje."
$found3 = $rng3.Find.Execute(
    $needle3, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $para3 = $rng3.Paragraphs(1)

    $noMark3 = $para3.Range
    $noMark3.MoveEnd(1, -1) | Out-Null

    $strikeRange3 = $d.Range($noMark3.Start, $noMark3.End - 1)
    $strikeRange3.Font.StrikeThrough = 1
}

Write-Output "Done"
